$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("164:164").Insert()

$ws.Cells.Item(164, 1).Value = 9
$ws.Cells.Item(164, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(164, 3).Value = "Metropolitana"
$ws.Cells.Item(164, 4).Value = 44824
$ws.Cells.Item(164, 5).Value = 13
$ws.Cells.Item(164, 6).Value = 100112026
$ws.Cells.Item(164, 7).Value = "Haba"
$ws.Cells.Item(164, 8).Value = "Sin especificar"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 70
$ws.Cells.Item(164, 11).Value = 14000
$ws.Cells.Item(164, 12).Value = 14000
$ws.Cells.Item(164, 13).Value = 14000
$ws.Cells.Item(164, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(164, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(164, 16).Value = 560
$ws.Cells.Item(164, 17).Value = 25
$ws.Cells.Item(164, 18).Value = "Hortaliza"

Write-Output "done"
